$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '245.19'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-0.74%'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '2.84%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.107'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.49%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05713'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.99%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.498'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.12%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8192'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.75%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8534'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '1.02%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0005978'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.31%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1330'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-0.70%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06935'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.78%'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-0.86%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09403'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.16%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001511'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.05%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.04033'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-13.37%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006215'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.04%'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-2.52%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.008'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.22%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.319'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '12.80%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '1.07%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '0.39%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.564'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-4.82%'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '1.70%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001216'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-2.60%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004471'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-2.58%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.00009895'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '3.12%'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '3.64%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03725'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.005917'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '73.06%'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-22.17%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002299'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-13.53%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009544'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '6.65%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005135'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-3.10%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000750'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.03%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1010'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-8.21%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002506'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-4.32%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002099'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.03%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001999'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.03%'
